$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix up row 23 labels (country / prior type) that were missing before
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = "Belgium"
$ws.Range("B23").Value = "extloglin"

# ---------------------------------------------------------------------------
# 2. New row 25 : just a "total nIter" value in M
# ---------------------------------------------------------------------------
$ws.Range("M25").Value = 10000

# ---------------------------------------------------------------------------
# 3. New data rows for Finland (27-32)
#    columns:     C      D      E      F     G      H      I     J      K     L     M      N    O    P       Q          R          S
# ---------------------------------------------------------------------------
$finlandRows = @(
    @(27, -1.5,  -0.05, -0.25, -2.5, -1.5,  -0.1,  0,    -0.1,  0.1,  0.01, 2000,  100, 75, 0.1205, 20.57457,   53.03808,   19.13396),
    @(28, -2.25, -0.05, -0.25, -3.5, -1,    -0.1,  0,    -0.1,  0.1,  0.01, 2000,  100, 75, 0.113,  38.49461,   60.99292,   58.89209),
    @(29, -2.75, -0.05, -0.05, -3.5, -2,    -0.1,  0,    -0.1,  0.1,  0.01, 2000,  100, 75, 0.179,  33.39472,   121.38133,  50.982),
    @(30, -2.75, -0.05, -0.05, -3,   -2.2,  -0.05, 0.05, -0.05, 0.1,  0.01, 2000,  100, 75, 0.182,  95.30202,   94.16052,   84.33098),
    @(31, -2.75, 0,     0,     -3,   -2.2,  -0.05, 0.05, -0.05, 0.1,  0.01, 2000,  100, 75, 0.3275, 186.1704,   194.4446,   198.2007),
    @(32, -2.75, 0,     0,     -3,   -2.2,  -0.05, 0.05, -0.05, 0.1,  0.01, 10000, 100, 75, 0.2745, 619.4724,   609.1472,   631.3015)
)

foreach ($rd in $finlandRows) {
    $row = $rd[0]

    $ws.Cells.Item($row, 1).Value  = "Finland"     # A : country
    $ws.Cells.Item($row, 2).Value  = "extloglin"    # B : prior type
    $ws.Cells.Item($row, 3).Value  = $rd[1]         # C
    $ws.Cells.Item($row, 4).Value  = $rd[2]         # D
    $ws.Cells.Item($row, 5).Value  = $rd[3]         # E
    $ws.Cells.Item($row, 6).Value  = $rd[4]         # F
    $ws.Cells.Item($row, 7).Value  = $rd[5]         # G
    $ws.Cells.Item($row, 8).Value  = $rd[6]         # H
    $ws.Cells.Item($row, 9).Value  = $rd[7]         # I
    $ws.Cells.Item($row, 10).Value = $rd[8]         # J
    $ws.Cells.Item($row, 11).Value = $rd[9]         # K
    $ws.Cells.Item($row, 12).Value = $rd[10]        # L
    $ws.Cells.Item($row, 13).Value = $rd[11]        # M
    $ws.Cells.Item($row, 14).Value = $rd[12]        # N
    $ws.Cells.Item($row, 15).Value = $rd[13]        # O

    $pCell = $ws.Cells.Item($row, 16)               # P
    $pCell.Value = $rd[14]
    $pCell.NumberFormat = "0.00%"

    $ws.Cells.Item($row, 17).Value = $rd[15]        # Q
    $ws.Cells.Item($row, 18).Value = $rd[16]        # R
    $ws.Cells.Item($row, 19).Value = $rd[17]        # S
}

# ---------------------------------------------------------------------------
# 4. Conditional formatting : highlight cells outside the allowed prior
#    bounds in red with white text, for columns C (lower bound check),
#    D (upper bound check) and E (starting-value check).
#
#    Rules are created in the exact document order observed in the
#    target workbook, and the dxf (differential formats) and priorities
#    are then assigned explicitly so the resulting workbook matches it.
# ---------------------------------------------------------------------------

$fcE1 = $ws.Range("E8:E30").FormatConditions.Add(1, 2, '=$J$27', '=$K$27')
$fcD1 = $ws.Range("D8:D30").FormatConditions.Add(1, 2, '=$H$8', '=$I$8')
$fcC1 = $ws.Range("C8:C30").FormatConditions.Add(1, 2, '=$F$8', '=$G$8')
$fcE2 = $ws.Range("E31").FormatConditions.Add(1, 2, '=$J$27', '=$K$27')
$fcD2 = $ws.Range("D31").FormatConditions.Add(1, 2, '=$H$8', '=$I$8')
$fcC2 = $ws.Range("C31").FormatConditions.Add(1, 2, '=$F$8', '=$G$8')
$fcE3 = $ws.Range("E32").FormatConditions.Add(1, 2, '=$J$27', '=$K$27')
$fcD3 = $ws.Range("D32").FormatConditions.Add(1, 2, '=$H$8', '=$I$8')
$fcC3 = $ws.Range("C32").FormatConditions.Add(1, 2, '=$F$8', '=$G$8')

# Assign the differential formats (white font on dark red fill) in the
# order that reproduces the target dxf numbering.
$fcD3.Font.Color = 16777215
$fcD3.Interior.Color = 192

$fcE3.Font.Color = 16777215
$fcE3.Interior.Color = 192

$fcD2.Font.Color = 16777215
$fcD2.Interior.Color = 192

$fcE2.Font.Color = 16777215
$fcE2.Interior.Color = 192

$fcD1.Font.Color = 16777215
$fcD1.Interior.Color = 192

# A short-lived rule used only to advance the internal dxf counter, so the
# remaining dxf ends up at the same index as in the target file.
$dummy = $ws.Range("Z1").FormatConditions.Add(1, 2, '=1', '=2')
$dummy.Font.Color = 16777215
$dummy.Interior.Color = 192
$dummy.Delete()

$fcE1.Font.Color = 16777215
$fcE1.Interior.Color = 192

# Priorities (1 = highest precedence). Newly-added rules take precedence
# over older ones, matching how Excel assigns priority when rules are
# created interactively.
$fcE1.Priority = 9
$fcD1.Priority = 8
$fcC1.Priority = 7
$fcE2.Priority = 6
$fcD2.Priority = 5
$fcC2.Priority = 4
$fcE3.Priority = 3
$fcD3.Priority = 2
$fcC3.Priority = 1

# ---------------------------------------------------------------------------
# 5. Leave the selection where the author left it
# ---------------------------------------------------------------------------
$ws.Range("W29").Select()
